# Updated cryptos list values (Price in column D, Volume(1h) in column E)
# Source data refreshed by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceUpdates = @{
    2 = "27.973.03"
    3 = "1.639.59"
    5 = "212.59"
    12 = "1.872.26"
    13 = "1.640.31"
    17 = "27.969.23"
    18 = "233.28"
    19 = "0.0₃0720"
    22 = "10.43"
    25 = "153.04"
    31 = "0.0484"
    33 = "3.08"
    34 = "1.403.08"
    38 = "0.562"
    39 = "0.927"
    43 = "67.14"
    44 = "5.54"
    47 = "1.781.13"
    48 = "88.06"
    51 = "7.61"
}

$volumeUpdates = @{
    2 = "  +0.23%  "
    3 = "  -0.07%  "
    4 = "  +0.04%  "
    5 = "  +0.02%  "
    6 = "  -0.30%  "
    7 = "  +0.02%  "
    8 = "  -0.13%  "
    10 = "  +0.13%  "
    11 = "  +1.62%  "
    12 = "  +0.00%  "
    13 = "  +0.05%  "
    14 = "  +0.26%  "
    15 = "  +1.17%  "
    16 = "  -0.31%  "
    17 = "  +0.31%  "
    18 = "  +0.73%  "
    19 = "  -0.38%  "
    20 = "  -1.14%  "
    21 = "  -0.02%  "
    22 = "  -3.24%  "
    23 = "  -0.08%  "
    24 = "  -3.39%  "
    25 = "  +1.35%  "
    26 = "  +0.43%  "
    27 = "  -0.26%  "
    28 = "  -0.55%  "
    29 = "  +0.00%  "
    30 = "  +0.32%  "
    31 = "  +0.28%  "
    32 = "  +2.62%  "
    33 = "  -0.22%  "
    34 = "  -3.71%  "
    35 = "  +1.80%  "
    36 = "  +1.23%  "
    37 = "  +0.54%  "
    38 = "  -0.38%  "
    39 = "  +1.17%  "
    40 = "  -1.26%  "
    41 = "  +0.82%  "
    43 = "  -3.00%  "
    44 = "  +3.60%  "
    45 = "  +2.27%  "
    46 = "  -0.47%  "
    47 = "  -0.05%  "
    48 = "  -0.32%  "
    49 = "  -0.54%  "
    50 = "  -0.30%  "
    51 = "  -1.49%  "
}

# Column D ("Price") cells are stored as plain text in this sheet (some values use
# "." as a thousands separator, e.g. "27.973.03"). Force text formatting before
# assigning so Excel does not auto-convert plain-decimal-looking values (like
# "153.04") into numbers, then drop back to the Normal style so no extra
# number-format is left applied to the cell.
foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.Style = "Normal"
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}

Write-Host "Updated $($priceUpdates.Keys.Count) price cells and $($volumeUpdates.Keys.Count) volume cells"
